# Generate Report for handback
# Updates the status of the "5fe91086-255c-4682-9774-a18afabf88ac.md" file
# (row 3) from "Not yet handed off" to "Handed back" across the Overview,
# zh-cn, and de-de sheets, and records new handback datetimes.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 5fe91086... file
$wsOverview.Range("B3").Value = "Handed back"
$wsOverview.Range("C3").Value = "Handed back"

# zh-cn sheet: row 3 is the 5fe91086... file
$wsZhCn.Range("B3").Value = "Handed back"
$wsZhCn.Range("G3").Value = "2016-01-08 14:18:47"

# de-de sheet: row 3 is the 5fe91086... file
$wsDeDe.Range("B3").Value = "Handed back"
$wsDeDe.Range("G3").Value = "2016-01-08 14:19:09"
